$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Terminé ?") holds boolean cells driven by linked Form-Control
# CheckBoxes (xl/ctrlProps/*.xml -> fmlaLink). Checking a checkbox in Excel
# sets its linked cell to TRUE. Mark the following "to do" rows as done by
# ticking the corresponding checkboxes (E7, E9, E10, E11, E17-E21).
$doneCells = @("E7", "E9", "E10", "E11", "E17", "E18", "E19", "E20", "E21")

foreach ($cellRef in $doneCells) {
    $ws.Range($cellRef).Value = $true
}
